# Actualización 11 de Mayo - Tarde
# Swap/correct the Materia (E) and Docente (F) values for a handful of
# rows on the "Blancos" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blancos")

# Row 28
$ws.Range("E28").Value = "APLICA ESTRUCTURAS DE DATOS CON UN LENGUAJE DE PROGRAMACIÓN"
$ws.Range("F28").Value = "De Jesús Orduña Sofía del Pilar"

# Row 29
$ws.Range("E29").Value = "CONSTRUYE ALGORITMOS PARA LA SOLUCIÓN DE PROBLEMAS"
$ws.Range("F29").Value = "Acevedo Rendón Ismael Arturo"

# Row 30
$ws.Range("E30").Value = "APLICA ESTRUCTURAS DE CONTROL CON UN LENGUAJE DE PROGRAMACIÓN"
$ws.Range("F30").Value = "Acevedo Rendón Ismael Arturo"

# Row 33
$ws.Range("E33").Value = "LECTURA, EXPRESIÓN ORAL Y ESCRITA II"
$ws.Range("F33").Value = "Medina Tolentino Francisco"

# Row 34
$ws.Range("E34").Value = "GEOMETRÍA Y TRIGONOMETRÍA"
$ws.Range("F34").Value = "Santiago Hernández Mariana"

# Row 35
$ws.Range("E35").Value = "GEOMETRÍA Y TRIGONOMETRÍA"
$ws.Range("F35").Value = "Santiago Hernández Mariana"

# Row 36
$ws.Range("E36").Value = "LECTURA, EXPRESIÓN ORAL Y ESCRITA II"
$ws.Range("F36").Value = "Medina Tolentino Francisco"
